$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices with locale-specific
# grouping dots, trailing zeros that matter, etc). Force each target
# cell to Text format before writing so Excel keeps the literal string
# instead of silently parsing it into a floating point number.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.787.17'
$ws.Range('E2').Value = '  -0.85%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.758.43'
$ws.Range('E3').Value = '  -2.86%  '

$ws.Range('E4').Value = '  -0.46%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.83'
$ws.Range('E5').Value = '  -2.39%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4232'
$ws.Range('E7').Value = '  -4.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3624'
$ws.Range('E8').Value = '  -2.65%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.41'
$ws.Range('E9').Value = '  -5.19%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07471'
$ws.Range('E10').Value = '  -2.90%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.084'
$ws.Range('E11').Value = '  -2.93%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.40%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.68'
$ws.Range('E13').Value = '  -5.48%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.063'
$ws.Range('E14').Value = '  -3.65%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.290'
$ws.Range('E15').Value = '  -2.76%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.759.09'
$ws.Range('E16').Value = '  -3.89%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.96'
$ws.Range('E17').Value = '  -2.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001054'
$ws.Range('E18').Value = '  -2.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06365'
$ws.Range('E19').Value = '  -1.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9997'
$ws.Range('E20').Value = '  -0.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.98'
$ws.Range('E21').Value = '  -2.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.928'
$ws.Range('E22').Value = '  -5.38%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.796.09'
$ws.Range('E23').Value = '  -1.01%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.20'
$ws.Range('E24').Value = '  -4.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.099'
$ws.Range('E25').Value = '  +0.53%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.28'
$ws.Range('E26').Value = '  +1.55%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.20'
$ws.Range('E27').Value = '  -1.98%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952.30'
$ws.Range('E28').Value = '  -3.67%  '

$ws.Range('E29').Value = '  -8.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.90'
$ws.Range('E30').Value = '  -2.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.111'
$ws.Range('E31').Value = '  -6.98%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.682'
$ws.Range('E32').Value = '  +0.43%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.547'
$ws.Range('E33').Value = '  -5.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08848'
$ws.Range('E34').Value = '  -4.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.20'
$ws.Range('E35').Value = '  -6.77%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02288'
$ws.Range('E36').Value = '  -2.17%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06045'
$ws.Range('E37').Value = '  -2.34%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2094'
$ws.Range('E38').Value = '  -3.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6311'
$ws.Range('E39').Value = '  -3.75%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.949'
$ws.Range('E40').Value = '  -4.09%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.176'
$ws.Range('E41').Value = '  -1.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9991'
$ws.Range('E42').Value = '  -0.42%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.868'
$ws.Range('E43').Value = '  -2.80%  '

$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.398'
$ws.Range('E44').Value = '  +1.02%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('E45').Value = '  -5.38%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5857'
$ws.Range('E46').Value = '  -3.48%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.681'
$ws.Range('E47').Value = '  -2.31%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.978'
$ws.Range('E48').Value = '  -2.89%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.91'
$ws.Range('E49').Value = '  -3.04%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.177'
$ws.Range('E50').Value = '  +2.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06816'
$ws.Range('E51').Value = '  -2.18%  '
